$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header labels in row 1
$ws.Range("C1").Value = "TentCount"
$ws.Range("D1").Value = "StructureCount"
$ws.Range("E1").Value = "VehicleCount"

# Data for columns C (TentCount), D (StructureCount), E (VehicleCount)
# for existing rows 2-18 (years 2008-2024), plus a brand-new row 19 (year 2025).
$data = @(
    @{ Row = 2;  Year = 2008; Sweep = 29;   Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 3;  Year = 2009; Sweep = 21;   Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 4;  Year = 2010; Sweep = 44;   Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 5;  Year = 2011; Sweep = 68;   Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 6;  Year = 2012; Sweep = 86;   Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 7;  Year = 2013; Sweep = 193;  Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 8;  Year = 2014; Sweep = 226;  Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 9;  Year = 2015; Sweep = 198;  Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 10; Year = 2016; Sweep = 696;  Tent = 0;    Structure = 0;   Vehicle = 0 },
    @{ Row = 11; Year = 2017; Sweep = 500;  Tent = 1602; Structure = 400; Vehicle = 134 },
    @{ Row = 12; Year = 2018; Sweep = 484;  Tent = 2113; Structure = 365; Vehicle = 333 },
    @{ Row = 13; Year = 2019; Sweep = 977;  Tent = 2515; Structure = 467; Vehicle = 62 },
    @{ Row = 14; Year = 2020; Sweep = 251;  Tent = 426;  Structure = 81;  Vehicle = 1 },
    @{ Row = 15; Year = 2021; Sweep = 53;   Tent = 352;  Structure = 0;   Vehicle = 0 },
    @{ Row = 16; Year = 2022; Sweep = 924;  Tent = 843;  Structure = 0;   Vehicle = 0 },
    @{ Row = 17; Year = 2023; Sweep = 2205; Tent = 1638; Structure = 304; Vehicle = 193 },
    @{ Row = 18; Year = 2024; Sweep = 2504; Tent = 4147; Structure = 36;  Vehicle = 38 },
    @{ Row = 19; Year = 2025; Sweep = 1259; Tent = 2068; Structure = 0;   Vehicle = 2 }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Year
    $ws.Cells.Item($r.Row, 2).Value = $r.Sweep
    $ws.Cells.Item($r.Row, 3).Value = $r.Tent
    $ws.Cells.Item($r.Row, 4).Value = $r.Structure
    $ws.Cells.Item($r.Row, 5).Value = $r.Vehicle
}

# Autofit the three newly-added columns so their widths are sized to the content
$ws.Columns("C:E").AutoFit() | Out-Null
